$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column F (reuse the bold/bordered header style from E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# time_taken values for rows 2-24
$timestamps = @(
    "2021-10-05 13:38:51.420528",
    "2021-10-05 13:38:51.420538",
    "2021-10-05 13:38:51.420541",
    "2021-10-05 13:38:51.420543",
    "2021-10-05 13:38:51.420545",
    "2021-10-05 13:38:51.420547",
    "2021-10-05 13:38:51.420549",
    "2021-10-05 13:38:51.420551",
    "2021-10-05 13:38:51.420581",
    "2021-10-05 13:38:51.420583",
    "2021-10-05 13:38:51.420585",
    "2021-10-05 13:38:51.420587",
    "2021-10-05 13:38:51.420589",
    "2021-10-05 13:38:51.420591",
    "2021-10-05 13:38:51.420593",
    "2021-10-05 13:38:51.420594",
    "2021-10-05 13:38:51.420597",
    "2021-10-05 13:38:51.420599",
    "2021-10-05 13:38:51.420601",
    "2021-10-05 13:38:51.420603",
    "2021-10-05 13:38:51.420605",
    "2021-10-05 13:38:51.420607",
    "2021-10-05 13:38:51.420609"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
